$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New period labels (reversed order of the old "Periodo Mora" column)
$ws.Range("E16").Value = "1709"
$ws.Range("E17").Value = "1712"
$ws.Range("E18").Value = "1801"
$ws.Range("E19").Value = "1802"
$ws.Range("E20").Value = "1803"

# Swap the "Valor Mora" amounts to match the re-ordered periods
$ws.Range("F16").Value = 29509
$ws.Range("F20").Value = 19673
